# Commit: "Added Service Pricing Column"
# The second worksheet was renamed from its placeholder title
# "I don't know what to upload." to "Service Pricing".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("I don't know what to upload.")
$ws.Name = "Service Pricing"
